$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Explicit portrait orientation on the page setup.
$sec.PageSetup.Orientation = 0

# Add the default footer: a 3-column layout table (left / center / right)
# followed by an empty Footer-styled paragraph, mirroring the existing
# header1.xml table-less layout but as a table-based footer.
$ftr = $sec.Footers.Item(1)
$xml = @"
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="TableNormal"/><w:bidiVisual w:val="0"/><w:tblW w:w="0" w:type="auto"/><w:tblLayout w:type="fixed"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3120"/><w:gridCol w:w="3120"/><w:gridCol w:w="3120"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:ind w:left="-115"/><w:jc w:val="left"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:jc w:val="center"/></w:pPr></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3120" w:type="dxa"/><w:tcMar/></w:tcPr><w:p><w:pPr><w:pStyle w:val="Header"/><w:bidi w:val="0"/><w:ind w:right="-115"/><w:jc w:val="right"/></w:pPr></w:p></w:tc></w:tr></w:tbl><w:p><w:pPr><w:pStyle w:val="Footer"/><w:bidi w:val="0"/></w:pPr></w:p>
"@
$ftr.Range.InsertXML($xml)
